# Swap the names of the two worksheets ("Planilha1" <-> "Clientes").
# The underlying sheets (tab position / data) do not move - only the
# names attached to each tab are exchanged. A temporary name is needed
# so the two renames don't collide.
$wb = $excel.ActiveWorkbook

$sheetA = $wb.Worksheets.Item("Planilha1")   # currently 1st tab
$sheetB = $wb.Worksheets.Item("Clientes")    # currently 2nd tab

$tempName = "__tmp_swap__"
$sheetA.Name = $tempName
$sheetB.Name = "Planilha1"
$sheetA.Name = "Clientes"

# After the rename, the sheet that used to be called "Clientes" (2nd
# tab) is now named "Planilha1" - update its view: scroll position and
# selected cell.
$ws = $wb.Worksheets.Item("Planilha1")
$ws.Activate()
$ws.Range("K169").Select()
$excel.ActiveWindow.ScrollRow = 142
$excel.ActiveWindow.ScrollColumn = 1
